$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.902.43"
$ws.Range("E2").Value = "  -0.22%  "
$ws.Range("D3").Value = "1.549.43"
$ws.Range("E3").Value = "  -0.34%  "
$ws.Range("E4").Value = "  -0.36%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "206.11"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.33%  "
$ws.Range("E6").Value = "  +0.85%  "
$ws.Range("E7").Value = "  -0.36%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.04"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.72%  "
$ws.Range("E9").Value = "  -0.41%  "
$ws.Range("E10").Value = "  +0.77%  "
$ws.Range("E11").Value = "  -0.49%  "
$ws.Range("D12").Value = "1.770.15"
$ws.Range("E12").Value = "  -0.26%  "
$ws.Range("D13").Value = "1.554.16"
$ws.Range("E13").Value = "  +0.08%  "
$ws.Range("E14").Value = "  +0.91%  "
$ws.Range("E15").Value = "  +0.59%  "
$ws.Range("D16").Value = "26.902.08"
$ws.Range("E16").Value = "  -0.20%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.63"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.26%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "217.23"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.30%  "
$ws.Range("E19").Value = "  +2.41%  "
$ws.Range("E20").Value = "  +0.01%  "
$ws.Range("E21").Value = "  -0.37%  "
$ws.Range("E22").Value = "  +0.70%  "
$ws.Range("E23").Value = "  +0.39%  "
$ws.Range("E24").Value = "  -0.93%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.49"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.23%  "
$ws.Range("E26").Value = "  -0.30%  "
$ws.Range("E27").Value = "  +0.44%  "
$ws.Range("E28").Value = "  +0.61%  "
$ws.Range("E29").Value = "  -0.32%  "
$ws.Range("E30").Value = "  +1.72%  "
$ws.Range("E31").Value = "  -1.17%  "
$ws.Range("E32").Value = "  -0.29%  "
$ws.Range("E33").Value = "  +4.27%  "
$ws.Range("D34").Value = "1.411.62"
$ws.Range("E34").Value = "  +2.37%  "
$ws.Range("E35").Value = "  +2.44%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.965"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.91%  "
$ws.Range("E37").Value = "  -0.08%  "
$ws.Range("E38").Value = "  +0.42%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.528"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.96%  "
$ws.Range("E40").Value = "  -0.37%  "
$ws.Range("E41").Value = "  -0.38%  "
$ws.Range("E42").Value = "  +3.77%  "
$ws.Range("E43").Value = "  +0.63%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.29"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.39%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "64.50"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.92%  "
$ws.Range("E46").Value = "  -0.04%  "
$ws.Range("D47").Value = "1.684.31"
$ws.Range("E47").Value = "  -0.25%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.13"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.98%  "
$ws.Range("E49").Value = "  +1.50%  "
$ws.Range("E50").Value = "  +5.06%  "
$ws.Range("E51").Value = "  +0.49%  "
